$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 3.1
$ws.Range("I2").Value = 2.5
$ws.Range("J2").Value = 4
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("Q2").Value = 1.88
$ws.Range("R2").Value = 1.98
$ws.Range("S2").Value = 2.5
$ws.Range("T2").Value = 1.5
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 1.17
$ws.Range("Y2").Value = 1.57
$ws.Range("Z2").Value = 2.25
$ws.Range("AA2").Value = 2.05
$ws.Range("AB2").Value = 1.7
$ws.Range("AI2").Value = 6.5
$ws.Range("AJ2").Value = 5.5
$ws.Range("AL2").Value = 67
$ws.Range("AM2").Value = 501
$ws.Range("AN2").Value = 6.5
$ws.Range("AR2").Value = 23
$ws.Range("AS2").Value = 41

# Row 3
$ws.Range("G3").Value = 2.8
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 2.15
$ws.Range("J3").Value = 3.25
$ws.Range("K3").Value = 2.3
$ws.Range("L3").Value = 2.67
$ws.Range("O3").Value = 1.19
$ws.Range("P3").Value = 3.75
$ws.Range("S3").Value = 1.57
$ws.Range("T3").Value = 2.1
$ws.Range("W3").Value = 2.37
$ws.Range("X3").Value = 1.45
$ws.Range("AA3").Value = 1.52
$ws.Range("AB3").Value = 2.2
$ws.Range("AC3").Value = 11.75
$ws.Range("AD3").Value = 16.5
$ws.Range("AE3").Value = 10.5
$ws.Range("AF3").Value = 32
$ws.Range("AG3").Value = 21
$ws.Range("AI3").Value = 14.5
$ws.Range("AJ3").Value = 7.4
$ws.Range("AL3").Value = 45
$ws.Range("AM3").Value = 250
$ws.Range("AN3").Value = 10.25
$ws.Range("AP3").Value = 9
$ws.Range("AQ3").Value = 21
$ws.Range("AR3").Value = 15.5
$ws.Range("AS3").Value = 22

# Row 5
$ws.Range("G5").Value = 1.67
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 4.5
$ws.Range("J5").Value = 2.38
$ws.Range("K5").Value = 2.1
$ws.Range("L5").Value = 5.5
$ws.Range("N5").Value = 10
$ws.Range("T5").Value = 1.75
$ws.Range("W5").Value = 3.75
$ws.Range("X5").Value = 1.25
$ws.Range("AC5").Value = 6.5
$ws.Range("AD5").Value = 7.5
$ws.Range("AF5").Value = 13
$ws.Range("AJ5").Value = 7
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 11
$ws.Range("AO5").Value = 23
$ws.Range("AP5").Value = 15
$ws.Range("AS5").Value = 41

